$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade record row (row 5) mirroring the existing rows' layout
$ws.Cells.Item(5, 1).Value = 10049.799999999999
$ws.Cells.Item(5, 2).Value = 9983.91
$ws.Cells.Item(5, 3).Value = 282.89999999999998
$ws.Cells.Item(5, 4).Value = 284.76
$ws.Cells.Item(5, 5).Value = $false
$ws.Cells.Item(5, 6).Value = 0.66
$ws.Cells.Item(5, 7).Value = 42609.506064814814
$ws.Cells.Item(5, 8).Value = $true

# Match the date-formatted style used in column G for the other rows
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)
